# Apply "Penalty Reward System" (unfinished) changes to the PO data workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

# Update B11 (requested quantity for week of 45123.99999999999) from 180 to 60
$ws1.Range("B11").Value = 60

# Remove row 12 entirely (it is deleted in the new version, shrinking the
# used range from A1:B12 down to A1:B11)
$ws1.Rows.Item(12).Delete()

# --- Sheet 2: "Monthly Trend" ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")

# Update B6 (requested quantity for month of 45138.99999999999) from 240 to 100
$ws2.Range("B6").Value = 100
